$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 148 (shifts existing rows 148-244 down to 149-245)
$ws.Range("A148").EntireRow.Insert()

# Populate the new row 148 with the new weekly record
$ws.Range("A148").Value = 5
$ws.Range("B148").Value = "Macroferia Regional de Talca"
$ws.Range("C148").Value = "Maule"
$ws.Range("D148").Value = 44704
$ws.Range("E148").Value = 7
$ws.Range("F148").Value = "Fruta"
$ws.Range("G148").Value = 100101
$ws.Range("H148").Value = "Berries"
$ws.Range("I148").Value = 100101007
$ws.Range("J148").Value = "Kiwi"
$ws.Range("K148").Value = "Hayward"
$ws.Range("L148").Value = "Primera"
$ws.Range("M148").Value = 180
$ws.Range("N148").Value = 10000
$ws.Range("O148").Value = 10000
$ws.Range("P148").Value = 10000
$ws.Range("Q148").Value = '$/bandeja 18 kilos'
$ws.Range("R148").Value = "Provincia de Curicó"
$ws.Range("S148").Value = 556
$ws.Range("T148").Value = 18
